# edit.ps1 -- PowerShell-style PowerPoint COM-interop script
#
# Reproduces the target edit:
#   1. Three tables (on the slides that hold the "Table_0" custom style)
#      get re-pointed from the deck's custom table style
#      {2BC82CB9-8142-4E1E-A91D-AEB278FBCE4E} to the built-in style
#      {784901C7-8794-49E5-AC2F-36F04FD0779C}.
#   2. The presentation's theme colour scheme (the "Integral" / Red Violet
#      theme used by the slide master) is swapped for the theme colours
#      that used to live in the deck's other theme part (the stock
#      "Office Theme" colours) -- i.e. the two embedded themes trade
#      their colour palettes.
#
# NOTE: PowerPoint's object model does not expose a generic "set table
# style" assignment -- the host explicitly requires Table.ApplyStyle(id).
# Likewise theme colours are only reachable/settable one swatch at a
# time via Master.ColorScheme.Colors(i).RGB (there is no supported way
# to rename the colour/theme, so only the RGB swatches are updated).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the three tables that use the custom "Table_0" style.
# ---------------------------------------------------------------------
$oldStyleId = "{2BC82CB9-8142-4E1E-A91D-AEB278FBCE4E}"
$newStyleId = "{784901C7-8794-49E5-AC2F-36F04FD0779C}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            $tbl.ApplyStyle($newStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colour palette used by the slide master / deck.
#    Order of ColorScheme.Colors(i) follows the OOXML clrScheme layout:
#    1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
#    Values below are the "Office Theme" palette that the deck's other
#    theme part carried before the swap.
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

$colorScheme.Colors(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1
$colorScheme.Colors(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1
$colorScheme.Colors(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2
$colorScheme.Colors(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2
$colorScheme.Colors(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1
$colorScheme.Colors(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2
$colorScheme.Colors(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3
$colorScheme.Colors(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4
$colorScheme.Colors(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5
$colorScheme.Colors(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6
$colorScheme.Colors(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink
$colorScheme.Colors(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink
